$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 3 and 4 (US/Dataset changed, task/status unchanged) ---
$ws.Range("A3").Value = "2022-09-06 13:43:03"
$ws.Range("B3").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C3").Value = "MER ATS"

$ws.Range("A4").Value = "2022-09-06 13:43:04"
$ws.Range("B4").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C4").Value = "MER ATS"

# --- Update row 5 (US/Dataset changed, task text changed to HTS TST, status now "processando...") ---
$ws.Range("A5").Value = "2022-09-06 13:43:10"
$ws.Range("B5").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C5").Value = "MER ATS"
$ws.Range("D5").Value = "Buscar valores para cada indicador: DSD HTS TST"
$ws.Range("E5").Value = "processando..."

# --- Append new row 6 ---
$ws.Range("A6").Value = "2022-09-06 13:44:00"
$ws.Range("B6").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C6").Value = "MER ATS"
$ws.Range("D6").Value = "Buscar valores para cada indicador: DSD HTS INDEX"
$ws.Range("E6").Value = "ok"

# --- Append new row 7 ---
$ws.Range("A7").Value = "2022-09-06 13:44:29"
$ws.Range("B7").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C7").Value = "MER ATS"
$ws.Range("D7").Value = "Buscar valores para cada indicador: DSD HTS SELF"
$ws.Range("E7").Value = "ok"

# --- Append new row 8 ---
$ws.Range("A8").Value = "2022-09-06 13:44:41"
$ws.Range("B8").Value = "MER_ATS_Xipamanine_12"
$ws.Range("C8").Value = "MER ATS"
$ws.Range("D8").Value = "Buscar valores para cada indicador: DSD TB STAT"
$ws.Range("E8").Value = "ok"
